# "set all eu files to download" — mark every G2 (EU) file row as download=1,
# add the autofilter over the G2 table, and switch the active tab from G2 to G1.

$wb = $excel.ActiveWorkbook
$wsG2 = $wb.Worksheets.Item("G2")

# Set download = 1 for every data row in column B (header already present in B1;
# B26 already had a value of 1, everything else was blank).
$rows = 2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28
foreach ($r in $rows) {
    $wsG2.Cells.Item($r, 2).Value = 1
}

# Turn on the AutoFilter for the G2 table.
$wsG2.Range("A1:L28").AutoFilter() | Out-Null

# The AutoFilter range is tracked by a hidden workbook-level, sheet-scoped
# defined name in the saved file.
$fdb = $wsG2.Names.Add("_xlnm._FilterDatabase", "='G2'!`$A`$1:`$L`$28")
$fdb.Visible = $false

# Update G2's remembered selection to B2:B28 (still selected, just not the
# active tab anymore once G1 is activated below).
$wsG2.Range("B2:B28").Select() | Out-Null

# Make G1 the active sheet/tab (was G2).
$wb.Worksheets.Item("G1").Activate()
